# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 9 (pushing existing rows 9-40 down to
# 10-41) and populate it with the new observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 9..40 down by one, creating a blank row 9.
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with the latest weekly record.
$ws.Range("A9").Value = 9
$ws.Range("B9").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44742
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112035
$ws.Range("G9").Value = "Bruselas (repollito)"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 52
$ws.Range("K9").Value = 20000
$ws.Range("L9").Value = 22000
$ws.Range("M9").Value = 21000
$ws.Range("N9").Value = '$/malla 15 kilos'
$ws.Range("O9").Value = "Hijuelas"
$ws.Range("P9").Value = 1400
$ws.Range("Q9").Value = 15
$ws.Range("R9").Value = "Hortaliza"
